# Updated smt results and added figures
# Applies the "SMT" section changes + the new "no SMT" data rows (31-35) to
# the first worksheet ("data-superpg") of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data-superpg")
$ws.Activate()

# ---------------------------------------------------------------------
# Section 1 (rows 1-6): "SMT" block.
# Row 3 ("4 hyperthreads/2cores") is merged into row 2's annotation text
# and removed; row 2 now reads "averaged over 100 iterations x 4
# hyperthreads/2cores". Row 4 keeps its row number/content (the node
# command), row 6 headers stay the same text.
# ---------------------------------------------------------------------
$ws.Range("A2").Value2 = "averaged over 100 iterations x 4 hyperthreads/2cores"
$ws.Range("A3").ClearContents()

# ---------------------------------------------------------------------
# Section 2 (rows 21-35): "no SMT" block.
# Row 22 becomes "averaged over 100 iterations x 2cores"; row 23
# ("2 processes/2cores") is removed. Rows 27-30 are relabeled/updated
# (usr-only counters), and four new rows (31-34) plus the relocated
# "elapse time" summary row (35) are appended.
# ---------------------------------------------------------------------
$ws.Range("A22").Value2 = "averaged over 100 iterations x 2cores"
$ws.Range("A23").ClearContents()

# --- Row 27: CPU_CLK_UNHALTED.THREAD_P,usr ---
$ws.Range("A27").Value2 = "CPU_CLK_UNHALTED.THREAD_P,usr"
$ws.Range("A27").Style = "Normal"
$ws.Range("B27").Value2 = 18272473320.9249
$ws.Range("C27").Value2 = 18266890799.77
$ws.Range("D27").Value2 = 18080593038.599899
$ws.Range("E27").Value2 = 18050384645.215
$ws.Range("G27").Formula = "=C27/B27"
$ws.Range("H27").Formula = "=D27/B27"
$ws.Range("I27").Formula = "=E27/B27"

# --- Row 28: DTLB_LOAD_MISSES.WALK_PENDING,usr ---
$ws.Range("A28").Value2 = "DTLB_LOAD_MISSES.WALK_PENDING,usr"
$ws.Range("A28").Style = "Normal"
$ws.Range("B28").Value2 = 87875652.090000004
$ws.Range("C28").Value2 = 88256127.049999997
$ws.Range("D28").Value2 = 82528229.590000004
$ws.Range("E28").Value2 = 75380325.400000006
$ws.Range("G28").Formula = "=C28/B28"
$ws.Range("H28").Formula = "=D28/B28"
$ws.Range("I28").Formula = "=E28/B28"

# --- Row 29: DTLB_STORE_MISSES.WALK_PENDING,usr ---
$ws.Range("A29").Value2 = "DTLB_STORE_MISSES.WALK_PENDING,usr"
$ws.Range("A29").Style = "Normal"
$ws.Range("B29").Value2 = 70347981.894999996
$ws.Range("C29").Value2 = 70363091.879999995
$ws.Range("D29").Value2 = 70266266.909999996
$ws.Range("E29").Value2 = 67740858.879999995
$ws.Range("G29").Formula = "=C29/B29"
$ws.Range("H29").Formula = "=D29/B29"
$ws.Range("I29").Formula = "=E29/B29"

# --- Row 30: ITLB_MISSES.WALK_PENDING,usr ---
$ws.Range("A30").Value2 = "ITLB_MISSES.WALK_PENDING,usr"
$ws.Range("A30").Style = "Normal"
$ws.Range("B30").Value2 = 82042447.620000005
$ws.Range("C30").Value2 = 74864658.760000005
$ws.Range("D30").Value2 = 63275600.880000003
$ws.Range("E30").Value2 = 66214448.664999999
$ws.Range("G30").Formula = "=C30/B30"
$ws.Range("H30").Formula = "=D30/B30"
$ws.Range("I30").Formula = "=E30/B30"

# --- Row 31 (new): ICACHE_64B.IFTAG_STALL,usr ---
$ws.Range("A31").Value2 = "ICACHE_64B.IFTAG_STALL,usr"
$ws.Range("B31").Value2 = 339995070.05000001
$ws.Range("C31").Value2 = 331442039.17000002
$ws.Range("D31").Value2 = 128115512.36
$ws.Range("E31").Value2 = 128899346.38
$ws.Range("G31").Formula = "=C31/B31"
$ws.Range("H31").Formula = "=D31/B31"
$ws.Range("I31").Formula = "=E31/B31"

# --- Row 32 (new): CPU_CLK_UNHALTED.THREAD_P ---
$ws.Range("A32").Value2 = "CPU_CLK_UNHALTED.THREAD_P"
$ws.Range("B32").Value2 = 18585172163.605
$ws.Range("C32").Value2 = 18590811046.455002
$ws.Range("D32").Value2 = 18396180939.18
$ws.Range("E32").Value2 = 18377656471.400002
$ws.Range("G32").Formula = "=C32/B32"
$ws.Range("H32").Formula = "=D32/B32"
$ws.Range("I32").Formula = "=E32/B32"

# --- Row 33 (new): INST_RETIRED.ANY_P ---
$ws.Range("A33").Value2 = "INST_RETIRED.ANY_P"
$ws.Range("B33").Value2 = 46793515450.989899
$ws.Range("C33").Value2 = 46813589276.580002
$ws.Range("D33").Value2 = 46803825874.809898
$ws.Range("E33").Value2 = 46807495795.360001
$ws.Range("G33").Formula = "=C33/B33"
$ws.Range("H33").Formula = "=D33/B33"
$ws.Range("I33").Formula = "=E33/B33"

# --- Row 34 (new): INST_RETIRED.ANY_P,usr ---
$ws.Range("A34").Value2 = "INST_RETIRED.ANY_P,usr"
$ws.Range("B34").Value2 = 46153563034.5299
$ws.Range("C34").Value2 = 46176789322.934898
$ws.Range("D34").Value2 = 46171617375.480003
$ws.Range("E34").Value2 = 46175753079.730003
$ws.Range("G34").Formula = "=C34/B34"
$ws.Range("H34").Formula = "=D34/B34"
$ws.Range("I34").Formula = "=E34/B34"

# --- Row 35 (new): elapse time (relocated summary row, red font like old A30) ---
$ws.Range("A35").Value2 = "elapse time"
$ws.Range("A35").Style = $ws.Range("A19").Style
$ws.Range("B35").Value2 = 531.21600000000001
$ws.Range("C35").Value2 = 531.245
$ws.Range("D35").Value2 = 525.77350000000001
$ws.Range("E35").Value2 = 525.07050000000004
$ws.Range("G35").Formula = "=C35/B35"
$ws.Range("H35").Formula = "=D35/B35"
$ws.Range("I35").Formula = "=E35/B35"

# ---------------------------------------------------------------------
# View: scroll position + active selection, matching the author's final
# on-screen state after adding the figures.
# ---------------------------------------------------------------------
$ws.Range("A8").Select()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K33").Select()

$wb.Save()
